$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 106, shifting existing rows 106:117 down to 107:118
$ws.Rows.Item(106).Insert()

# Populate the new row 106 with data
$ws.Cells.Item(106, 1).Value = 1
$ws.Cells.Item(106, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(106, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(106, 4).Value = 45223
$ws.Cells.Item(106, 5).Value = 15
$ws.Cells.Item(106, 6).Value = 100112012
$ws.Cells.Item(106, 7).Value = "Espinaca"
$ws.Cells.Item(106, 8).Value = "Sin especificar"
$ws.Cells.Item(106, 9).Value = "Primera"
$ws.Cells.Item(106, 10).Value = 275
$ws.Cells.Item(106, 11).Value = 1400
$ws.Cells.Item(106, 12).Value = 1500
$ws.Cells.Item(106, 13).Value = 1445
$ws.Cells.Item(106, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(106, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(106, 16).Value = 482
$ws.Cells.Item(106, 17).Value = 3
$ws.Cells.Item(106, 18).Value = "Hortaliza"
